# Fruta / hortaliza, semanal
# This edit permutes the weekly data across rows 2-15: the values in
# columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) are reshuffled among
# the rows, while the rest of each row (market/region/category metadata)
# stays fixed per row. We read the current values for all rows first,
# then write them back out in the permuted order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a block per source row.
$cols = @("D", "J", "K", "L", "M", "P")

# Capture current ("before") values for rows 2..15 so we can freely
# reassign without clobbering source data mid-update.
$before = @{}
for ($r = 2; $r -le 15; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# Mapping: destination row -> source row whose D/J/K/L/M/P values it
# should receive (derived from matching unique Fecha (D) values
# between the before and after states of the workbook).
$mapping = @{
    2  = 12
    3  = 2
    4  = 15
    5  = 6
    6  = 14
    7  = 13
    8  = 4
    9  = 7
    10 = 3
    11 = 9
    12 = 10
    13 = 8
    14 = 11
    15 = 5
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $before[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
